$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Cover sheet: row heights, label text, and logo picture resize/reposition
# ---------------------------------------------------------------------------
$cover = $wb.Worksheets.Item("Cover")

# Row height changes
$cover.Rows.Item(3).RowHeight = 20
$cover.Rows.Item(13).RowHeight = 10
$cover.Rows.Item(15).RowHeight = 40

# Label text change
$cover.Range("B9").Value = "Customer Name:"

# Remove the 2nd and 3rd pictures (Image 2 / Image 3), keep Image 1
$cover.Shapes.Item(3).Delete()
$cover.Shapes.Item(2).Delete()

# Move + resize remaining logo picture (Image 1) to row 15 (1-based) / col B
$logo = $cover.Shapes.Item(1)
$logo.Top = $cover.Cells.Item(15, 2).Top
$logo.Left = $cover.Cells.Item(15, 2).Left
$logo.Width = 90
$logo.Height = 22.5

# ---------------------------------------------------------------------------
# LOE sheet: refresh the "Notes" (column K) copy for every task row
# ---------------------------------------------------------------------------
$loe = $wb.Worksheets.Item("LOE")

$notes = @{
    2  = "Includes stakeholder interviews and requirements documentation deliverable"
    3  = "Customer must provide access to current infrastructure for evaluation"
    4  = "Identifies compliance gaps and remediation requirements"
    5  = "Customer must provide API documentation and integration endpoints"
    6  = "Deliverable: Security design document with controls and compliance mapping"
    7  = "Deliverable: Data model diagrams and migration strategy"
    8  = "Includes environment provisioning and configuration baseline"
    9  = "Largest effort item - varies significantly based on feature complexity"
    10 = "Includes schema design stored procedures and initial data seeding"
    11 = "Customer integration endpoints must be available for testing"
    12 = "Includes responsive design for desktop and mobile devices"
    13 = "Deliverable: Functional authentication authorization and audit logging"
    14 = "Automated test suite for regression testing"
    15 = "Requires access to customer test environments and systems"
    16 = "Load testing with expected concurrent user volumes"
    17 = "Third-party security assessment may require additional budget"
    18 = "Customer business users must be available for UAT period"
    19 = "Final checks before production deployment"
    20 = "Includes deployment runbook and rollback procedures"
    21 = "Customer responsible for data quality and validation"
    22 = "Includes daily monitoring and rapid issue response during hypercare"
    23 = "Includes system administration and troubleshooting procedures"
    24 = "Role-based training sessions - customer provides training space"
    25 = "Deliverable: User guides administrator guides and runbooks"
    26 = "Comprehensive handover to customer operations team"
    27 = "Establishes SLA metrics and monitoring thresholds"
    28 = "Solution architecture review and technical escalation support"
    29 = "Includes weekly status reporting and risk management"
}

foreach ($row in $notes.Keys) {
    $loe.Range("K$row").Value = $notes[$row]
}
